$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Fix typo in existing shared string: remove stray space in
# "仁寶電腦工業股份有限公 司•" -> "仁寶電腦工業股份有限公司•"
$ws.Range("B3").Value = "仁寶電腦工業股份有限公司•"

# Insert a new "property_category" column before the existing "date"
# column (H), shifting date/legislator_name/legislator_id one column
# to the right (H->I, I->J, J->K).
$ws.Columns.Item(8).Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H11").Value = "stock"
